$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new log entry (Post 40 - Semaphores and Counting Semaphores).
# Cells are written in this specific order (F, C, E, B, D) so that the
# workbook's shared-string table grows in the same sequence as the source
# edit: dev.to link, title, hashnode link.
$ws.Range("F50").Value = "https://dev.to/rahulmishra05/semaphores-and-counting-semaphores-operating-system-m03-p05-3fbo"
$ws.Range("C50").Value = "Semaphores and Counting Semaphores | Operating System - M03 P05"
$ws.Range("E50").Value = "https://programmingport.hashnode.dev/semaphores-and-counting-semaphores-or-operating-system-m03-p05"
$ws.Range("B50").Value = 40
$ws.Range("D50").Value = 44168

# Grow the worksheet table ("Table2") so it covers the new row; this also
# keeps the AutoFilter range and the sheet dimension in sync.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("B10:F50"))

# Match the author's final selection/scroll position after the edit.
$ws.Range("E50").Select() | Out-Null
